$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 96: Date, Game, ModCount - appended after row 95
# Force the date column to be stored as literal text (matching the
# existing rows) rather than letting Excel auto-convert it to a date
# serial number.
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = "2026/02/14"
$ws.Range("A96").Style = "Normal"

$ws.Range("B96").Value = "逃离鸭科夫"
$ws.Range("C96").Value = 1202

# Match the centered alignment style used by the rest of the data rows
$ws.Range("A96:C96").HorizontalAlignment = -4108
$ws.Range("A96:C96").VerticalAlignment = -4108
